$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new log row (row 11) as plain text in every column, mirroring
# the existing rows (all stored as text, even the numeric-looking ones).
# A leading apostrophe forces text interpretation for the cells that would
# otherwise be auto-detected as a number (empty "A11" / numeric "C11");
# resetting the style back to Normal afterwards drops the quote-prefix
# formatting flag while the cell keeps its text content/type.
$ws.Range("A11").Value = "'"
$ws.Range("A11").Style = "Normal"

$ws.Range("B11").Value = "يامن"

$ws.Range("C11").Value = "'23"
$ws.Range("C11").Style = "Normal"

$ws.Range("D11").Value = "الجزائري"
$ws.Range("E11").Value = "الرحلة 2"
$ws.Range("F11").Value = "C3"
$ws.Range("G11").Value = "NRC"
$ws.Range("H11").Value = "٠٢‏/٠٥‏/٢٠٢٥ ٠٢:١٣:٤٩ م"
